# Auto-applying cell value updates to refresh market-price derived columns
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets, per scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 79.07143000000001
$ws.Range("I9").Value = 58.916668
$ws.Range("K9").Value = 58.916668
$ws.Range("M9").Value = 110.083332
$ws.Range("H21").Value = 363906.34
$ws.Range("I21").Value = 1700
$ws.Range("J21").Value = 545009.5
$ws.Range("K21").Value = 1700
$ws.Range("L21").Value = 545009.5
$ws.Range("M21").Value = -1232
$ws.Range("N21").Value = -545945.5
$ws.Range("H23").Value = 363906.34
$ws.Range("I23").Value = 1700
$ws.Range("J23").Value = 545009.5
$ws.Range("K23").Value = 1700
$ws.Range("L23").Value = 545009.5
$ws.Range("M23").Value = -1466
$ws.Range("N23").Value = -545477.5
$ws.Range("H29").Value = 1035.75
$ws.Range("I29").Value = 71.5
$ws.Range("K29").Value = 214.5
$ws.Range("M29").Value = 66.5
$ws.Range("H38").Value = 522.8333
$ws.Range("I38").Value = 96.28570999999999
$ws.Range("J38").Value = 1120
$ws.Range("K38").Value = 288.85713
$ws.Range("L38").Value = 3360
$ws.Range("M38").Value = 83.14287000000002
$ws.Range("N38").Value = -4104
$ws.Range("H58").Value = 6877.778
$ws.Range("J58").Value = 9833.333000000001
$ws.Range("L58").Value = 29499.999
$ws.Range("N58").Value = -29799.999
$ws.Range("H87").Value = 25142.857
$ws.Range("J87").Value = 26000
$ws.Range("L87").Value = 26000
$ws.Range("N87").Value = -28496
$ws.Range("H90").Value = 25142.857
$ws.Range("J90").Value = 26000
$ws.Range("L90").Value = 78000
$ws.Range("N90").Value = -90480
$ws.Range("H92").Value = 982.72
$ws.Range("I92").Value = 917.65
$ws.Range("K92").Value = 917.65
$ws.Range("M92").Value = 330.35
$ws.Range("H103").Value = 566.3333
$ws.Range("I103").Value = 500
$ws.Range("J103").Value = 599.5
$ws.Range("K103").Value = 1500
$ws.Range("L103").Value = 1798.5
$ws.Range("M103").Value = -914
$ws.Range("N103").Value = -2970.5
$ws.Range("H111").Value = 1277.7778
$ws.Range("I111").Value = 2200
$ws.Range("J111").Value = 1162.5
$ws.Range("K111").Value = 6600
$ws.Range("L111").Value = 3487.5
$ws.Range("M111").Value = -3533
$ws.Range("N111").Value = -9621.5
$ws.Range("H112").Value = 1304.3334
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 1340.3572
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 4021.0716
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -6237.071599999999
$ws.Range("H115").Value = 342.5
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1047.7059
$ws.Range("I2").Value = 805.5
$ws.Range("J2").Value = 1393.7142
$ws.Range("K2").Value = 805.5
$ws.Range("L2").Value = 1393.7142
$ws.Range("M2").Value = -692.5
$ws.Range("N2").Value = -1619.7142
$ws.Range("H32").Value = 3831.9814
$ws.Range("I32").Value = 3357.8164
$ws.Range("K32").Value = 3357.8164
$ws.Range("M32").Value = -3070.8164
$ws.Range("H74").Value = 1942.826
$ws.Range("I74").Value = 1610.2778
$ws.Range("K74").Value = 1610.2778
$ws.Range("M74").Value = -736.2778000000001
$ws.Range("H77").Value = 1942.826
$ws.Range("I77").Value = 1610.2778
$ws.Range("K77").Value = 8051.389
$ws.Range("M77").Value = -3683.389
$ws.Range("H102").Value = 1780
$ws.Range("I102").Value = 1780
$ws.Range("K102").Value = 1780
$ws.Range("M102").Value = -158
$ws.Range("H116").Value = 1047.7059
$ws.Range("I116").Value = 805.5
$ws.Range("J116").Value = 1393.7142
$ws.Range("K116").Value = 805.5
$ws.Range("L116").Value = 1393.7142
$ws.Range("M116").Value = 1488.5
$ws.Range("N116").Value = -5981.7142
$ws.Range("H122").Value = 1669.6364
$ws.Range("I122").Value = 1228.6666
$ws.Range("J122").Value = 2198.8
$ws.Range("K122").Value = 3685.9998
$ws.Range("L122").Value = 6596.400000000001
$ws.Range("M122").Value = -1235.9998
$ws.Range("N122").Value = -11496.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1047.7059
$ws.Range("I3").Value = 805.5
$ws.Range("J3").Value = 1393.7142
$ws.Range("K3").Value = 805.5
$ws.Range("L3").Value = 1393.7142
$ws.Range("M3").Value = -691.5
$ws.Range("N3").Value = -1621.7142
$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -53744
$ws.Range("H99").Value = 749.625
$ws.Range("I99").Value = 582.8333
$ws.Range("K99").Value = 582.8333
$ws.Range("M99").Value = 915.1667
$ws.Range("H105").Value = 2201.8
$ws.Range("J105").Value = 2500
$ws.Range("L105").Value = 2500
$ws.Range("N105").Value = -5994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 691
$ws.Range("I16").Value = 629.3333
$ws.Range("J16").Value = 765
$ws.Range("K16").Value = 629.3333
$ws.Range("L16").Value = 765
$ws.Range("M16").Value = -342.3333
$ws.Range("N16").Value = -1339
$ws.Range("H31").Value = 6063271
$ws.Range("I31").Value = 2793.926
$ws.Range("J31").Value = 33335416
$ws.Range("K31").Value = 2793.926
$ws.Range("L31").Value = 33335416
$ws.Range("M31").Value = -2498.926
$ws.Range("N31").Value = -33336006
$ws.Range("H34").Value = 6063271
$ws.Range("I34").Value = 2793.926
$ws.Range("J34").Value = 33335416
$ws.Range("K34").Value = 2793.926
$ws.Range("L34").Value = 33335416
$ws.Range("M34").Value = -2591.926
$ws.Range("N34").Value = -33335820
$ws.Range("H113").Value = 691
$ws.Range("I113").Value = 629.3333
$ws.Range("J113").Value = 765
$ws.Range("K113").Value = 629.3333
$ws.Range("L113").Value = 765
$ws.Range("M113").Value = 1540.6667
$ws.Range("N113").Value = -5105

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3888.7368
$ws.Range("I3").Value = 2218.6
$ws.Range("K3").Value = 6655.799999999999
$ws.Range("M3").Value = -6543.799999999999
$ws.Range("H64").Value = 1000
$ws.Range("I64").Value = 500
$ws.Range("J64").Value = 1500
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 4500
$ws.Range("M64").Value = -1230
$ws.Range("N64").Value = -5040
$ws.Range("H67").Value = 1000
$ws.Range("I67").Value = 500
$ws.Range("J67").Value = 1500
$ws.Range("K67").Value = 1500
$ws.Range("L67").Value = 4500
$ws.Range("M67").Value = -564
$ws.Range("N67").Value = -6372
$ws.Range("H80").Value = 2260
$ws.Range("J80").Value = 2260
$ws.Range("L80").Value = 6780
$ws.Range("N80").Value = -8652
$ws.Range("H83").Value = 2260
$ws.Range("J83").Value = 2260
$ws.Range("L83").Value = 20340
$ws.Range("N83").Value = -29700
$ws.Range("H114").Value = 2080.4783
$ws.Range("I114").Value = 1615.5714
$ws.Range("J114").Value = 2283.875
$ws.Range("K114").Value = 4846.7142
$ws.Range("L114").Value = 6851.625
$ws.Range("M114").Value = -1592.7142
$ws.Range("N114").Value = -13359.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 9666.666999999999
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 9666.666999999999
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 9666.666999999999
$ws.Range("M53").Value = ""
$ws.Range("N53").Value = -10928.667
$ws.Range("H57").Value = 19733.334
$ws.Range("J57").Value = 19733.334
$ws.Range("L57").Value = 19733.334
$ws.Range("N57").Value = -21373.334
$ws.Range("H102").Value = 1008.875
$ws.Range("I102").Value = 909.2
$ws.Range("J102").Value = 1175
$ws.Range("K102").Value = 909.2
$ws.Range("L102").Value = 1175
$ws.Range("M102").Value = 712.8
$ws.Range("N102").Value = -4419
$ws.Range("H132").Value = 79358.12
$ws.Range("I132").Value = 107048.42
$ws.Range("K132").Value = 321145.26
$ws.Range("M132").Value = -318615.26

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 54040
$ws.Range("J75").Value = 54040
$ws.Range("L75").Value = 54040
$ws.Range("N75").Value = -55912
$ws.Range("H78").Value = 54040
$ws.Range("J78").Value = 54040
$ws.Range("L78").Value = 162120
$ws.Range("N78").Value = -171480
$ws.Range("H100").Value = 1272.2222
$ws.Range("I100").Value = 1243.75
$ws.Range("K100").Value = 1243.75
$ws.Range("M100").Value = -702.75
